# Add data for 2022-01-02
# Updates the "through December 24" date label to "through December 25" and
# refreshes the carjacking counts affected by the newly added day of data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet and update the column header label to reflect the new date.
$ws.Name = "Through 2021-12-25"
$ws.Range("B1").Value = "December 2021 (through December 25)"

# Englewood
$ws.Range("N3").Value = 3

# North Lawndale
$ws.Range("BV4").Value = 3

# Roseland
$ws.Range("AX13").Value = 3
$ws.Range("BJ13").Value = 8

# United Center
$ws.Range("Z17").Value = 1

# South Chicago
$ws.Range("B30").Value = 4

# Brighton Park
$ws.Range("B39").Value = 3

# Washington Park
$ws.Range("N47").Value = 1
$ws.Range("BV47").Value = 1

# Hyde Park
$ws.Range("BJ53").Value = 1

# Avondale
$ws.Range("BJ66").Value = 1

# Belmont Cragin
$ws.Range("Z67").Value = 2

# Hegewisch
$ws.Range("BJ78").Value = 1

# Logan Square
$ws.Range("AL82").Value = 1
$ws.Range("BJ82").Value = 2

# River North
$ws.Range("AX93").Value = 1
